$d = $word.ActiveDocument

# 1) "in C# or " -> "in C#, Ruby, or " (adds Ruby to the language list, right before
#    the spell-checked "Javascript" run, which we leave untouched).
$null = $d.Content.Find.Execute(
    "in C# or ", $true, $false, $false, $false, $false, $true, 1, $false,
    "in C#, Ruby, or ", 2)

# 2) Insert the new "From automation to YAGNI..." sentence between the first and
#    second sentence of the paragraph, and turn the following comma into a period
#    ("...and projects, with..." -> "...and projects. With...").
$null = $d.Content.Find.Execute(
    "our craft. These are lessons learned from past applications and projects, with an influence",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "our craft. From automation to YAGNI and TDD to collective ownership, these are lessons learned from past applications and projects. With an influence",
    2)

# 3) "many mentors, that will help" -> "many mentors, these thoughts will help"
$null = $d.Content.Find.Execute(
    "many mentors, that will help", $true, $false, $false, $false, $false, $true, 1, $false,
    "many mentors, these thoughts will help", 2)
